# Update "Homicidios (Ni cuanto commit)" values.
# The author retyped several values in columns B (Número de Homicidios) and
# D (Homicidios Hombres) that had been stored as pseudo-decimal numbers
# (Spanish thousands-separator typed as a decimal point, e.g. 1.389 meaning
# 1,389) into plain integers (e.g. 1389).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2"  = 1389;  "D2"  = 1236
    "B3"  = 1172;  "D3"  = 103
    "B5"  = 1026
    "B6"  = 1068
    "B7"  = 1216;  "D7"  = 1025
    "B8"  = 1616;  "D8"  = 1446
    "B9"  = 181;   "D9"  = 1623
    "B10" = 1978
    "B11" = 2268;  "D11" = 2067
    "B12" = 2513;  "D12" = 228
    "B13" = 2474;  "D13" = 2264
    "B14" = 2375;  "D14" = 2195
    "B15" = 2394;  "D15" = 2207
    "B16" = 2111;  "D16" = 1936
    "B17" = 2144;  "D17" = 1964
    "B18" = 2394;  "D18" = 218
    "B19" = 2239
    "B20" = 1902
    "B21" = 1279
    "B22" = 1095
    "B23" = 121
    "B24" = 1066
    "B25" = 1121
    "B26" = 1038
    "B27" = 1032
    "B28" = 1007
    "B29" = 1005
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Reflect the selection recorded in the saved file: the user had just
# finished editing the B2:D32 block with B2 as the active cell.
$ws.Range("B2:D32").Select()
